# Scheduled market-data refresh: update computed price/profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 918.375
$ws.Range("I9").Value = 453.6
$ws.Range("J9").Value = 1693
$ws.Range("K9").Value = 453.6
$ws.Range("L9").Value = 1693
$ws.Range("M9").Value = -284.6
$ws.Range("N9").Value = -2031
$ws.Range("H33").Value = 13890221
$ws.Range("I33").Value = 22727644
$ws.Range("J33").Value = 2842.1428
$ws.Range("K33").Value = 22727644
$ws.Range("L33").Value = 2842.1428
$ws.Range("M33").Value = -22727415
$ws.Range("N33").Value = -3300.1428
$ws.Range("H69").Value = 4443
$ws.Range("I69").Value = 4443
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 13329
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -12455
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 4443
$ws.Range("I72").Value = 4443
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 39987
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -35619
$ws.Range("N72").ClearContents()
$ws.Range("H86").Value = 11844.223
$ws.Range("I86").Value = 10524.5
$ws.Range("K86").Value = 10524.5
$ws.Range("M86").Value = -9401.5
$ws.Range("H88").Value = 1333.7
$ws.Range("I88").Value = 1388.8
$ws.Range("J88").Value = 1278.6
$ws.Range("K88").Value = 1388.8
$ws.Range("L88").Value = 1278.6
$ws.Range("M88").Value = -982.8
$ws.Range("N88").Value = -2090.6
$ws.Range("H89").Value = 11844.223
$ws.Range("I89").Value = 10524.5
$ws.Range("K89").Value = 52622.5
$ws.Range("M89").Value = -47006.5
$ws.Range("H91").Value = 1333.7
$ws.Range("I91").Value = 1388.8
$ws.Range("J91").Value = 1278.6
$ws.Range("K91").Value = 1388.8
$ws.Range("L91").Value = 1278.6
$ws.Range("M91").Value = 15.20000000000005
$ws.Range("N91").Value = -4086.6
$ws.Range("H138").Value = 4299.7886
$ws.Range("J138").Value = 4510.2915
$ws.Range("L138").Value = 13530.8745
$ws.Range("N138").Value = -23810.8745

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2325
$ws.Range("I122").Value = 2205.7693
$ws.Range("K122").Value = 6617.3079
$ws.Range("M122").Value = -4167.3079
$ws.Range("H135").Value = 295214.5
$ws.Range("J135").Value = 295214.5
$ws.Range("L135").Value = 295214.5
$ws.Range("N135").Value = -305354.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 11449.889
$ws.Range("I22").Value = 14671.429
$ws.Range("J22").Value = 174.5
$ws.Range("K22").Value = 14671.429
$ws.Range("L22").Value = 174.5
$ws.Range("M22").Value = -14498.429
$ws.Range("N22").Value = -520.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1529.5
$ws.Range("I16").Value = 1661.3334
$ws.Range("J16").Value = 1134
$ws.Range("K16").Value = 1661.3334
$ws.Range("L16").Value = 1134
$ws.Range("M16").Value = -1374.3334
$ws.Range("N16").Value = -1708
$ws.Range("H105").Value = 9325.333000000001
$ws.Range("I105").Value = 10790.4
$ws.Range("K105").Value = 10790.4
$ws.Range("M105").Value = -9043.4
$ws.Range("H107").Value = 1474
$ws.Range("I107").Value = 1285.0714
$ws.Range("J107").Value = 1851.8572
$ws.Range("K107").Value = 1285.0714
$ws.Range("L107").Value = 1851.8572
$ws.Range("M107").Value = 634.9286
$ws.Range("N107").Value = -5691.8572
$ws.Range("H113").Value = 1529.5
$ws.Range("I113").Value = 1661.3334
$ws.Range("J113").Value = 1134
$ws.Range("K113").Value = 1661.3334
$ws.Range("L113").Value = 1134
$ws.Range("M113").Value = 508.6666
$ws.Range("N113").Value = -5474

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 9005
$ws.Range("I59").Value = 9005
$ws.Range("K59").Value = 27015
$ws.Range("M59").Value = -26475
$ws.Range("H60").Value = 5067
$ws.Range("I60").Value = 5067
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 15201
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -14950
$ws.Range("N60").ClearContents()
$ws.Range("H131").Value = 1527.12
$ws.Range("I131").Value = 1466.6666
$ws.Range("K131").Value = 4399.9998
$ws.Range("M131").Value = 640.0002000000004
$ws.Range("H139").Value = 18460.334
$ws.Range("I139").Value = 19517.875
$ws.Range("J139").Value = 10000
$ws.Range("K139").Value = 58553.625
$ws.Range("L139").Value = 30000
$ws.Range("M139").Value = -53413.625
$ws.Range("N139").Value = -40280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2970.4285
$ws.Range("I80").Value = 5598.3335
$ws.Range("J80").Value = 999.5
$ws.Range("K80").Value = 5598.3335
$ws.Range("L80").Value = 999.5
$ws.Range("M80").Value = -4600.3335
$ws.Range("N80").Value = -2995.5
$ws.Range("H83").Value = 2970.4285
$ws.Range("I83").Value = 5598.3335
$ws.Range("J83").Value = 999.5
$ws.Range("K83").Value = 27991.6675
$ws.Range("L83").Value = 4997.5
$ws.Range("M83").Value = -22999.6675
$ws.Range("N83").Value = -14981.5
$ws.Range("H132").Value = 761554.9399999999
$ws.Range("I132").Value = 5448.6
$ws.Range("K132").Value = 16345.8
$ws.Range("M132").Value = -13815.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 11099.272
$ws.Range("I68").Value = 17683
$ws.Range("K68").Value = 17683
$ws.Range("M68").Value = -16934
$ws.Range("H71").Value = 11099.272
$ws.Range("I71").Value = 17683
$ws.Range("K71").Value = 88415
$ws.Range("M71").Value = -84671
$ws.Range("H122").Value = 9894.706
$ws.Range("J122").Value = 8550
$ws.Range("L122").Value = 25650
$ws.Range("N122").Value = -30550
$ws.Range("H136").Value = 145863.62
$ws.Range("I136").Value = 20255.455
$ws.Range("K136").Value = 60766.36500000001
$ws.Range("M136").Value = -58216.36500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1655.9333
$ws.Range("I81").Value = 1712.7273
$ws.Range("J81").Value = 1499.75
$ws.Range("K81").Value = 3425.4546
$ws.Range("L81").Value = 2999.5
$ws.Range("M81").Value = -2364.4546
$ws.Range("N81").Value = -5121.5
$ws.Range("H84").Value = 1655.9333
$ws.Range("I84").Value = 1712.7273
$ws.Range("J84").Value = 1499.75
$ws.Range("K84").Value = 17127.273
$ws.Range("L84").Value = 14997.5
$ws.Range("M84").Value = -11823.273
$ws.Range("N84").Value = -25605.5
$ws.Range("H107").Value = 990.2
$ws.Range("I107").Value = 1036.5385
$ws.Range("K107").Value = 3109.6155
$ws.Range("M107").Value = -1189.6155
$ws.Range("H122").Value = 6999.909
$ws.Range("I122").Value = 2750
$ws.Range("J122").Value = 7944.3335
$ws.Range("K122").Value = 8250
$ws.Range("L122").Value = 23833.0005
$ws.Range("M122").Value = -5800
$ws.Range("N122").Value = -28733.0005
$ws.Range("H132").Value = 990314.0600000001
$ws.Range("I132").Value = 2154.4443
$ws.Range("K132").Value = 6463.3329
$ws.Range("M132").Value = -3933.3329

